$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed cryptocurrency price / 1h-volume data.
#
# All of the source cells hold plain text (e.g. "4.90", "0.999",
# "  +0.94%  "), even when that text happens to look like a number.
# If we just assign .Value, Excel auto-detects the numeric-looking
# strings and stores them as real numbers, which silently destroys
# meaningful trailing zeros (e.g. "4.90" would become 4.9). To avoid
# that, each cell is temporarily switched to a text number format
# before the value is written, then its original Style object is
# restored so formatting stays exactly as it was.

$cellStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.164.36'
$ws.Range('D2').Style = $cellStyle
$cellStyle = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('E2').Style = $cellStyle
$cellStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.381.52'
$ws.Range('D3').Style = $cellStyle
$cellStyle = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E3').Style = $cellStyle
$cellStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = $cellStyle
$cellStyle = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').Style = $cellStyle
$cellStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.54'
$ws.Range('D5').Style = $cellStyle
$cellStyle = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('E5').Style = $cellStyle
$cellStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.67'
$ws.Range('D6').Style = $cellStyle
$cellStyle = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('E6').Style = $cellStyle
$cellStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').Style = $cellStyle
$cellStyle = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('E7').Style = $cellStyle
$cellStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.369.82'
$ws.Range('D8').Style = $cellStyle
$cellStyle = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('E8').Style = $cellStyle
$cellStyle = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E9').Style = $cellStyle
$cellStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.173'
$ws.Range('D10').Style = $cellStyle
$cellStyle = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.68%  '
$ws.Range('E10').Style = $cellStyle
$cellStyle = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('E11').Style = $cellStyle
$cellStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.58'
$ws.Range('D12').Style = $cellStyle
$cellStyle = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('E12').Style = $cellStyle
$cellStyle = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('E13').Style = $cellStyle
$cellStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.17'
$ws.Range('D14').Style = $cellStyle
$cellStyle = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('E14').Style = $cellStyle
$cellStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.919.19'
$ws.Range('D15').Style = $cellStyle
$cellStyle = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('E15').Style = $cellStyle
$cellStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.31'
$ws.Range('D16').Style = $cellStyle
$cellStyle = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('E16').Style = $cellStyle
$cellStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.376.29'
$ws.Range('D17').Style = $cellStyle
$cellStyle = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('E17').Style = $cellStyle
$cellStyle = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('E18').Style = $cellStyle
$cellStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.89'
$ws.Range('D19').Style = $cellStyle
$cellStyle = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('E19').Style = $cellStyle
$cellStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '64.926.98'
$ws.Range('D20').Style = $cellStyle
$cellStyle = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('E20').Style = $cellStyle
$cellStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('D21').Style = $cellStyle
$cellStyle = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('E21').Style = $cellStyle
$cellStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '454.37'
$ws.Range('D22').Style = $cellStyle
$cellStyle = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('E22').Style = $cellStyle
$cellStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.90'
$ws.Range('D23').Style = $cellStyle
$cellStyle = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('E23').Style = $cellStyle
$cellStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.35'
$ws.Range('D24').Style = $cellStyle
$cellStyle = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +7.82%  '
$ws.Range('E24').Style = $cellStyle
$cellStyle = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.20%  '
$ws.Range('E25').Style = $cellStyle
$cellStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.36'
$ws.Range('D26').Style = $cellStyle
$cellStyle = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.04%  '
$ws.Range('E26').Style = $cellStyle
$cellStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.89'
$ws.Range('D27').Style = $cellStyle
$cellStyle = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('E27').Style = $cellStyle
$cellStyle = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.45%  '
$ws.Range('E28').Style = $cellStyle
$cellStyle = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('E29').Style = $cellStyle
$cellStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.15'
$ws.Range('D30').Style = $cellStyle
$cellStyle = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.85%  '
$ws.Range('E30').Style = $cellStyle
$cellStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '63.10'
$ws.Range('D32').Style = $cellStyle
$cellStyle = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.48%  '
$ws.Range('E32').Style = $cellStyle
$cellStyle = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E33').Style = $cellStyle
$cellStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '578.49'
$ws.Range('D34').Style = $cellStyle
$cellStyle = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('E34').Style = $cellStyle
$cellStyle = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E35').Style = $cellStyle
$cellStyle = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E36').Style = $cellStyle
$cellStyle = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.96%  '
$ws.Range('E37').Style = $cellStyle
$cellStyle = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('E38').Style = $cellStyle
$cellStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.72'
$ws.Range('D39').Style = $cellStyle
$cellStyle = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E39').Style = $cellStyle
$cellStyle = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('E40').Style = $cellStyle
$cellStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0740'
$ws.Range('D41').Style = $cellStyle
$cellStyle = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E41').Style = $cellStyle
$cellStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.089.86'
$ws.Range('D42').Style = $cellStyle
$cellStyle = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('E42').Style = $cellStyle
$cellStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0416'
$ws.Range('D43').Style = $cellStyle
$cellStyle = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('E43').Style = $cellStyle
$cellStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.76'
$ws.Range('D44').Style = $cellStyle
$cellStyle = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('E44').Style = $cellStyle
$cellStyle = $ws.Range('B45').Style
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Stellar'
$ws.Range('B45').Style = $cellStyle
$cellStyle = $ws.Range('C45').Style
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C45').Style = $cellStyle
$cellStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.134'
$ws.Range('D45').Style = $cellStyle
$cellStyle = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('E45').Style = $cellStyle
$cellStyle = $ws.Range('B46').Style
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Fetch.AI'
$ws.Range('B46').Style = $cellStyle
$cellStyle = $ws.Range('C46').Style
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('C46').Style = $cellStyle
$cellStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.45'
$ws.Range('D46').Style = $cellStyle
$cellStyle = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('E46').Style = $cellStyle
$cellStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.16'
$ws.Range('D47').Style = $cellStyle
$cellStyle = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.94%  '
$ws.Range('E47').Style = $cellStyle
$cellStyle = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E48').Style = $cellStyle
$cellStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.49'
$ws.Range('D49').Style = $cellStyle
$cellStyle = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.51%  '
$ws.Range('E49').Style = $cellStyle
$cellStyle = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('E50').Style = $cellStyle
$cellStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.28'
$ws.Range('D51').Style = $cellStyle
$cellStyle = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.72%  '
$ws.Range('E51').Style = $cellStyle
